# Weekly price-list update: insert a new record row for the market
# "Feria Lagunitas de Puerto Montt" (Cilantro) at row 186, pushing the
# existing rows 186-248 down to 187-249.
#
# The new row duplicates the (now shifted) row below it for the columns
# that stay constant (market/region/category/unit/origin/classification)
# and then gets its own data for Fecha / Volumen / Precios.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 186, shifting old rows 186:248 to 187:249.
$ws.Rows("186:186").Insert()

# Duplicate the row that is now at 187 (the old row 186) into the new
# blank row 186, so all the constant columns (A,B,C,E,F,G,H,I,N,O,Q,R)
# are already populated correctly.
$ws.Rows("187:187").Copy()
$ws.Rows("186:186").PasteSpecial()

# Now overwrite the data-specific cells of the new row 186 with the
# new observation's values.
$ws.Cells.Item(186, 4).Value = 44627   # D186 Fecha
$ws.Cells.Item(186, 10).Value = 120    # J186 Volumen
$ws.Cells.Item(186, 11).Value = 15000  # K186 Precio minimo
$ws.Cells.Item(186, 12).Value = 15000  # L186 Precio maximo
$ws.Cells.Item(186, 13).Value = 15000  # M186 Precio promedio ponderado
$ws.Cells.Item(186, 16).Value = 417    # P186 Precio $/Kg
